# "fixed date on minutes"
# The header currently shows the placeholder text "TEMP, 2015"; replace it
# with the actual meeting date, "March 3, 2015".

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    # Primary (default) header for this section.
    $hdr = $sec.Headers.Item(1)
    if ($hdr.Exists) {
        $hdr.Range.Find.Execute("TEMP, 2015", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "March 3, 2015", 2)
    }
}
